# Apply the "swap sheet_name/Author columns" edit to the readme table,
# and refresh the analysis timestamp on the Project Information sheet.

$wb = $excel.ActiveWorkbook

# --- 1. readme sheet: swap columns D (Author) and E (sheet_name) ---
$readme = $wb.Worksheets.Item("readme")

# Swap header labels (row 1)
$dHeader = $readme.Range("D1").Value2
$eHeader = $readme.Range("E1").Value2
$readme.Range("D1").Value = $eHeader
$readme.Range("E1").Value = $dHeader

# Swap data values for each row (rows 2 through 12)
for ($row = 2; $row -le 12; $row++) {
    $dCell = $readme.Cells.Item($row, 4)
    $eCell = $readme.Cells.Item($row, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}

# --- 2. Project Information sheet: update "Date of Analysis" timestamp ---
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Range("B11").Value = "2022-03-03 16:23:21.729804"
